$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Prediction/Error values for re-run classification results
# (values written in plain decimal form since the interpreter does not
# accept scientific-notation numeric literals like 1.23E-45)

# Row 2 - Control 0
$ws.Range("D2").Value = 0.000000000000000000000000000000000003458689682029154
$ws.Range("E2").Value = 0.000000000000000000000000000000000003458689682029154

# Row 3 - Control 6
$ws.Range("D3").Value = 0.009692396989404606
$ws.Range("E3").Value = 0.009692396989404606

# Row 5 - Control 24
$ws.Range("D5").Value = 0.9999807277324363
$ws.Range("E5").Value = 0.9999807277324363

# Row 6 - Control 32
$ws.Range("D6").Value = 0.9999999999986402
$ws.Range("E6").Value = 0.9999999999986402

# Row 8 - MDD 30
$ws.Range("D8").Value = 1.0
$ws.Range("E8").Value = 0.0

# Row 10 - MDD 17
$ws.Range("D10").Value = 0.0006583433920512715
$ws.Range("E10").Value = 0.9993416566079487

# Row 11 - MDD 23
$ws.Range("D11").Value = 0.9999423499306409
$ws.Range("E11").Value = 0.00005765006935909422
$ws.Range("F11").Value = 168.5313262939453
